$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("par jours")
$ws2 = $wb.Worksheets.Item("par semaines")

# --- New shared strings must be introduced in this exact order so that
# --- they land at shared-string table indices 40..55, matching the target.
$ws1.Range("D21").Value = "Documenter le travail de la semaine et se préparer pour la  semaine 6 si nécessaire"
$ws1.Range("D17").Value = "Se familliariser avec les dataGridView"
$ws1.Range("D22").Value = "Implémenter le formulaire permettant de saisir les informations pour créer un nouveau TPI"
$ws1.Range("D23").Value = "Implémenter la vérification de l'intégrité des données saisies"
$ws1.Range("D26").Value = "Documenter le travail de la semaine et se préparer pour la  semaine7 si nécessaire"
$ws1.Range("D24").Value = "Implémenter l'insertion du TPI saisi dans la DB"
$ws1.Range("D25").Value = "Ajouter au formulaire une liste de technologies qui peuvent être reliées au TPI ET modifier la requête d'ajout en conséquence"
$ws1.Range("B37").Value = "?"
$ws1.Range("D34").Value = "Rendu du projet"
$ws1.Range("D33").Value = "Finalisation documentation"
$ws1.Range("D35").Value = "Préparation de la défense"
$ws1.Range("D32").Value = "Implémenter le bouton permettant de supprimmer un TPI"
$ws1.Range("D27").Value = "Implémenter un formulaire permettant de modifer les données d'un TPI existant"
$ws1.Range("D31").Value = "Documenter le travail de la semaine et se préparer pour la  semaine8 si nécessaire"
$ws1.Range("D29").Value = "Implémenter la vérification des données saisies"
$ws1.Range("D30").Value = "Implémenter la modification des données du TPI dans la DB"

# --- Cells that reuse already-existing shared strings.
$ws1.Range("D19").Value = "Implémenter le formulaire affichant les informations détaillées d'un TPI  "
$ws1.Range("D20").Value = "Implémenter le formulaire affichant les informations détaillées d'un TPI  "
$ws1.Range("D28").Value = "Implémenter un formulaire permettant de modifer les données d'un TPI existant"
$ws1.Range("D36").Value = "Préparation de la défense"

# --- Fill in the week 6/7/8 dates in column B.
$ws1.Range("B22").Value = 44999
$ws1.Range("B23").Value = 45000
$ws1.Range("B24").Value = 45000
$ws1.Range("B25").Value = 45001
$ws1.Range("B26").Value = 45001
$ws1.Range("B27").Value = 45006
$ws1.Range("B28").Value = 45007
$ws1.Range("B29").Value = 45007
$ws1.Range("B30").Value = 45008
$ws1.Range("B31").Value = 45008
$ws1.Range("B32").Value = 45013
$ws1.Range("B33").Value = 45014
$ws1.Range("B34").Value = 45014
$ws1.Range("B35").Value = 45015
$ws1.Range("B36").Value = 45015

# --- Rows whose wrapped text now spans two lines get a taller row height.
$ws1.Rows.Item(22).RowHeight = 30
$ws1.Rows.Item(25).RowHeight = 30

# --- View/selection state: "par semaines" keeps a selection but is no longer
# --- the active tab; "par jours" becomes active, scrolled down, with E12 selected.
$ws2.Range("B7").Select()
$ws1.Activate()
$ws1.Range("E12").Select()
